# Weekly roll-down of the price-history table:
#   - A brand-new week's observation is inserted as row 3 (pushing the
#     previous rows 3-17 down to rows 4-18).
#   - Columns that are constant for every observation (Mercado ID,
#     Mercado, Región, Codreg, Categoría ID, Categoría, Variedad,
#     Calidad, Unidad de comercialización, Origen, Kg o Unidades,
#     Clasificación) are copied from the row immediately below (the
#     former row 2/row 3), which already carries the correct values.
#   - The variable columns (Fecha, Volumen, Precio mínimo, Precio
#     máximo, Precio promedio ponderado, Precio $/Kg) get the new
#     week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing data rows 3:17 down to 4:18, leaving row 3 empty.
$ws.Rows.Item(3).Insert()

# Copy the constant columns from row 4 (the row that used to be row 3)
# so the new row matches the rest of the table.
$ws.Range("A3").Value = $ws.Range("A4").Value()
$ws.Range("B3").Value = $ws.Range("B4").Value()
$ws.Range("C3").Value = $ws.Range("C4").Value()
$ws.Range("E3").Value = $ws.Range("E4").Value()
$ws.Range("F3").Value = $ws.Range("F4").Value()
$ws.Range("G3").Value = $ws.Range("G4").Value()
$ws.Range("H3").Value = $ws.Range("H4").Value()
$ws.Range("I3").Value = $ws.Range("I4").Value()
$ws.Range("N3").Value = $ws.Range("N4").Value()
$ws.Range("O3").Value = $ws.Range("O4").Value()
$ws.Range("Q3").Value = $ws.Range("Q4").Value()
$ws.Range("R3").Value = $ws.Range("R4").Value()

# New week's figures.
$ws.Range("D3").Value = 44764
$ws.Range("J3").Value = 120
$ws.Range("K3").Value = 7000
$ws.Range("L3").Value = 8000
$ws.Range("M3").Value = 7500
$ws.Range("P3").Value = 125
